$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update every cell that currently shows "Ready for handoff" to "In Translation"
# (the status text was refreshed while regenerating the report).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# The status columns are now narrower, so re-fit their width to the new
# (shorter) "In Translation" content that replaced "Ready for handoff".
$wsOverview.Columns.Item(5).EntireColumn.AutoFit()
$wsOverview.Columns.Item(6).EntireColumn.AutoFit()
$wsZhCn.Columns.Item(3).EntireColumn.AutoFit()
$wsDeDe.Columns.Item(3).EntireColumn.AutoFit()

# Nudge to the narrower width used in the published report.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
